# Updates cryptos list with refreshed price/volume data (GitHub Actions run).
# For cells whose new text looks like a plain number (e.g. "1.001"), the
# column's NumberFormat is forced to text ("@") before assigning the value so
# Excel does not silently reinterpret it as a numeric/date value, then
# ClearFormats() removes the temporary format again so the cell keeps the
# workbook's original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.937.60'
$ws.Range("E2").Value = '  -1.32%  '

$ws.Range("D3").Value = '1.776.19'
$ws.Range("E3").Value = '  -1.54%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.47'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5352'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.68%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3738'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.61%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07425'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.60'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.087'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.98%  '

$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.38'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.58%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.056'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.40%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.782.16'
$ws.Range("E15").Value = '  -1.06%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.199'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.70'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001049'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06404'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.71%  '

$ws.Range("E20").Value = '  +0.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.20'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.872'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.82%  '

$ws.Range("D23").Value = '27.942.88'
$ws.Range("E23").Value = '  -1.27%  '

$ws.Range("E24").Value = '  -3.08%  '

$ws.Range("E25").Value = '  -3.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.06'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.17'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.40%  '

$ws.Range("D28").Value = '1.971.53'
$ws.Range("E28").Value = '  -1.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.271'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.37%  '

$ws.Range("E30").Value = '  -3.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.108'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1040'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.91%  '

$ws.Range("E33").Value = '  -0.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.495'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2231'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06344'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02259'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.948'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.369'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6105'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.425'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.95'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.169'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.02%  '

$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.18'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.653'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.88%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5725'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.42'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.180'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.917'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06787'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.76%  '
